# Auto-generated from verified diff parse: apply all cell value updates
# across the 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 61
$ws.Range("H61").Value = 813.5
$ws.Range("I61").Value = 825.6667
$ws.Range("J61").Value = 777
$ws.Range("K61").Value = 2477.0001
$ws.Range("L61").Value = 2331
$ws.Range("M61").Value = -2305.0001
$ws.Range("N61").Value = -2675
# row 88
$ws.Range("H88").Value = 3149.4375
$ws.Range("I88").Value = 4865.375
$ws.Range("K88").Value = 4865.375
$ws.Range("M88").Value = -4459.375
# row 91
$ws.Range("H91").Value = 3149.4375
$ws.Range("I91").Value = 4865.375
$ws.Range("K91").Value = 4865.375
$ws.Range("M91").Value = -3461.375
# row 123
$ws.Range("H123").Value = 59992.5
$ws.Range("J123").Value = 59992.5
$ws.Range("L123").Value = 59992.5
$ws.Range("N123").Value = -69792.5

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 18882952
$ws.Range("I32").Value = 21752598
$ws.Range("K32").Value = 21752598
$ws.Range("M32").Value = -21752311
# row 61
$ws.Range("H61").Value = 6014.241
$ws.Range("I61").Value = 5074.4736
$ws.Range("K61").Value = 5074.4736
$ws.Range("M61").Value = -4862.4736
# row 74
$ws.Range("H74").Value = 1698.0476
$ws.Range("I74").Value = 1337.3334
$ws.Range("K74").Value = 1337.3334
$ws.Range("M74").Value = -463.3334
# row 77
$ws.Range("H77").Value = 1698.0476
$ws.Range("I77").Value = 1337.3334
$ws.Range("K77").Value = 6686.666999999999
$ws.Range("M77").Value = -2318.666999999999
# row 88
$ws.Range("H88").Value = 4909924
$ws.Range("J88").Value = 8335102
$ws.Range("L88").Value = 8335102
$ws.Range("N88").Value = -8335914
# row 91
$ws.Range("H91").Value = 4909924
$ws.Range("J91").Value = 8335102
$ws.Range("L91").Value = 8335102
$ws.Range("N91").Value = -8337910
# row 132
$ws.Range("H132").Value = 3136.375
$ws.Range("I132").Value = 3136.375
$ws.Range("K132").Value = 9409.125
$ws.Range("M132").Value = -6879.125
# row 136
$ws.Range("H136").Value = 6014.241
$ws.Range("I136").Value = 5074.4736
$ws.Range("K136").Value = 15223.4208
$ws.Range("M136").Value = -12673.4208

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 3547.7856
$ws.Range("I31").Value = 2383.1177
$ws.Range("J31").Value = 5347.727
$ws.Range("K31").Value = 2383.1177
$ws.Range("L31").Value = 5347.727
$ws.Range("M31").Value = -2088.1177
$ws.Range("N31").Value = -5937.727
# row 34
$ws.Range("H34").Value = 3547.7856
$ws.Range("I34").Value = 2383.1177
$ws.Range("J34").Value = 5347.727
$ws.Range("K34").Value = 2383.1177
$ws.Range("L34").Value = 5347.727
$ws.Range("M34").Value = -2181.1177
$ws.Range("N34").Value = -5751.727
# row 58
$ws.Range("H58").Value = 1523.7742
$ws.Range("I58").Value = 932.9545
$ws.Range("K58").Value = 932.9545
$ws.Range("M58").Value = -729.9545
# row 99
$ws.Range("H99").Value = 16035198
$ws.Range("I99").Value = 4073128.8
$ws.Range("J99").Value = 25006750
$ws.Range("K99").Value = 4073128.8
$ws.Range("L99").Value = 25006750
$ws.Range("M99").Value = -4071630.8
$ws.Range("N99").Value = -25009746
# row 126
$ws.Range("H126").Value = 16035198
$ws.Range("I126").Value = 4073128.8
$ws.Range("J126").Value = 25006750
$ws.Range("K126").Value = 12219386.4
$ws.Range("L126").Value = 75020250
$ws.Range("M126").Value = -12216916.4
$ws.Range("N126").Value = -75025190
# row 136
$ws.Range("H136").Value = 1523.7742
$ws.Range("I136").Value = 932.9545
$ws.Range("K136").Value = 2798.8635
$ws.Range("M136").Value = -248.8635000000004

$ws = $wb.Worksheets.Item("CUL")
# row 18
$ws.Range("H18").Value = 452.14285
$ws.Range("J18").Value = 795
$ws.Range("L18").Value = 2385
$ws.Range("N18").Value = -2723
# row 26
$ws.Range("H26").Value = 731.25
$ws.Range("J26").Value = 550
$ws.Range("L26").Value = 1650
$ws.Range("N26").Value = -2226
# row 32
$ws.Range("H32").Value = 200000000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 200000000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 600000000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -600000566
# row 46
$ws.Range("H46").Value = 335333.34
$ws.Range("I46").Value = 335333.34
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1006000.02
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1005909.02
$ws.Range("N46").ClearContents()
# row 112
$ws.Range("H112").Value = 1100
$ws.Range("I112").Value = 1100
$ws.Range("K112").Value = 3300
$ws.Range("M112").Value = -2192
# row 114
$ws.Range("H114").Value = 1202.75
$ws.Range("I114").Value = 487.25
$ws.Range("J114").Value = 1560.5
$ws.Range("K114").Value = 1461.75
$ws.Range("L114").Value = 4681.5
$ws.Range("M114").Value = 1792.25
$ws.Range("N114").Value = -11189.5
# row 124
$ws.Range("H124").Value = 2099.5
$ws.Range("I124").Value = 1466
$ws.Range("K124").Value = 4398
$ws.Range("M124").Value = 512

$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 40242.75
$ws.Range("I70").Value = 66669.5
$ws.Range("J70").Value = 6265.5
$ws.Range("K70").Value = 66669.5
$ws.Range("L70").Value = 6265.5
$ws.Range("M70").Value = -66399.5
$ws.Range("N70").Value = -6805.5
# row 73
$ws.Range("H73").Value = 40242.75
$ws.Range("I73").Value = 66669.5
$ws.Range("J73").Value = 6265.5
$ws.Range("K73").Value = 66669.5
$ws.Range("L73").Value = 6265.5
$ws.Range("M73").Value = -65733.5
$ws.Range("N73").Value = -8137.5
# row 80
$ws.Range("H80").Value = 14543680
$ws.Range("I80").Value = 64189.223
$ws.Range("J80").Value = 66669844
$ws.Range("K80").Value = 64189.223
$ws.Range("L80").Value = 66669844
$ws.Range("M80").Value = -63191.223
$ws.Range("N80").Value = -66671840
# row 83
$ws.Range("H83").Value = 14543680
$ws.Range("I83").Value = 64189.223
$ws.Range("J83").Value = 66669844
$ws.Range("K83").Value = 320946.115
$ws.Range("L83").Value = 333349220
$ws.Range("M83").Value = -315954.115
$ws.Range("N83").Value = -333359204
# row 104
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
# row 126
$ws.Range("H126").Value = 5938.6665
$ws.Range("I126").Value = 4984.4287
$ws.Range("K126").Value = 14953.2861
$ws.Range("M126").Value = -12483.2861
# row 134
$ws.Range("H134").Value = 91999.75
$ws.Range("J134").Value = 91999.75
$ws.Range("L134").Value = 275999.25
$ws.Range("N134").Value = -281069.25

$ws = $wb.Worksheets.Item("LTW")
# row 132
$ws.Range("H132").Value = 4743.943
$ws.Range("I132").Value = 4706.853
$ws.Range("J132").Value = 6005
$ws.Range("K132").Value = 14120.559
$ws.Range("L132").Value = 18015
$ws.Range("M132").Value = -11590.559
$ws.Range("N132").Value = -23075
# row 135
$ws.Range("H135").Value = 78545.45
$ws.Range("J135").Value = 78545.45
$ws.Range("L135").Value = 78545.45
$ws.Range("N135").Value = -88685.45
# row 136
$ws.Range("H136").Value = 5216.696
$ws.Range("I136").Value = 3090
$ws.Range("K136").Value = 9270
$ws.Range("M136").Value = -6720

$ws = $wb.Worksheets.Item("WVR")
# row 96
$ws.Range("H96").Value = 64925.25
$ws.Range("I96").Value = 85616.336
$ws.Range("J96").Value = 2852
$ws.Range("K96").Value = 85616.336
$ws.Range("L96").Value = 2852
$ws.Range("M96").Value = -84243.336
$ws.Range("N96").Value = -5598
# row 123
$ws.Range("H123").Value = 51818.09
$ws.Range("J123").Value = 51818.09
$ws.Range("L123").Value = 51818.09
$ws.Range("N123").Value = -61618.09
# row 125
$ws.Range("H125").Value = 69545.45
$ws.Range("J125").Value = 69545.45
$ws.Range("L125").Value = 69545.45
$ws.Range("N125").Value = -79385.45
